$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42604.891423611109
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"

$ws.Range("B5").Value = "Random"

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 57
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 59
$ws.Range("M5").Value = 41
